$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "001/TTT/AV1"
$ws.Range("C2").Value = "B219321"
$ws.Range("D2").Value = "JEMAA HORMI"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 2000
$ws.Range("J2").Value = 0
$ws.Range("O2").Value = 2000

# Row 3
$ws.Range("A3").Value = "001/TTT/AV1"
$ws.Range("C3").Value = "IR801997"
$ws.Range("D3").Value = "NOUBAIL MOHAMMED"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 2000
$ws.Range("J3").Value = 0
$ws.Range("O3").Value = 2000

# Row 4
$ws.Range("A4").Value = "001/TTT/AV1"
$ws.Range("C4").Value = "IB19558"
$ws.Range("D4").Value = "ZERNAKH ABDELLAH"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 2000
$ws.Range("J4").Value = 0
$ws.Range("O4").Value = 2000

# Row 5
$ws.Range("H5").Value = 6000
$ws.Range("J5").Value = 0
$ws.Range("O5").Value = 6000
